# Update countries & provincias Spain
# - "Asturias" moves up in the ranking, ahead of "Segovia" and "Leon"
#   (Segovia's old row becomes Leon's row+1, etc. -- i.e. the province
#   names in rows 26-28 shift while Asturias gets freshly updated counts)
# - Refresh the "last updated" timestamp string
# - Refresh a batch of per-province case/recovered/death counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "last updated" banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 11:52"

# --- Asturias jumps ahead of Segovia / Leon in the table -------------------
$ws.Range("A26").Value = "Asturias"
$ws.Range("A27").Value = "Segovia"
$ws.Range("A28").Value = "Leon"

# --- updated counts: Madrid ------------------------------------------------
$ws.Range("B4").Value = 59199
$ws.Range("C4").Value = 33032
$ws.Range("D4").Value = 18590
$ws.Range("E4").Value = 7577

# --- updated counts: Cataluna -----------------------------------------------
$ws.Range("B5").Value = 43802
$ws.Range("C5").Value = 15089
$ws.Range("D5").Value = 24466
$ws.Range("E5").Value = 4247

# --- updated counts: Castilla-La Mancha -------------------------------------
$ws.Range("B6").Value = 17321
$ws.Range("C6").Value = 4337
$ws.Range("D6").Value = 10844
$ws.Range("E6").Value = 2140

# --- updated counts: Castilla y Leon -----------------------------------------
$ws.Range("B7").Value = 16839
$ws.Range("C7").Value = 5614
$ws.Range("D7").Value = 9671
$ws.Range("E7").Value = 1554

# --- updated counts: Pais Vasco ---------------------------------------------
$ws.Range("B8").Value = 13044
$ws.Range("C8").Value = 7651
$ws.Range("D8").Value = 4269
$ws.Range("E8").Value = 1124

# --- updated counts: Andalucia -----------------------------------------------
$ws.Range("B9").Value = 11610
$ws.Range("C9").Value = 3569
$ws.Range("D9").Value = 6991
$ws.Range("E9").Value = 1050

# --- updated counts: Galicia --------------------------------------------------
$ws.Range("B10").Value = 8634
$ws.Range("C10").Value = 1625
$ws.Range("D10").Value = 6641
$ws.Range("E10").Value = 368

# --- updated counts: Aragon ----------------------------------------------------
$ws.Range("B14").Value = 5054
$ws.Range("C14").Value = 1530
$ws.Range("D14").Value = 2868
$ws.Range("E14").Value = 656

# --- updated counts: Navarra -----------------------------------------------------
$ws.Range("B15").Value = 4899
$ws.Range("C15").Value = 1316
$ws.Range("D15").Value = 3182
$ws.Range("E15").Value = 401

# --- updated counts: La Rioja -------------------------------------------------------
$ws.Range("B17").Value = 3792
$ws.Range("C17").Value = 1899
$ws.Range("D17").Value = 1565
$ws.Range("E17").Value = 298

# --- updated counts: Extremadura -----------------------------------------------------
$ws.Range("B22").Value = 3230
$ws.Range("C22").Value = 1051
$ws.Range("D22").Value = 1775
$ws.Range("E22").Value = 404

# --- updated counts: Asturias (now row 26) --------------------------------------------
$ws.Range("B26").Value = 2419
$ws.Range("C26").Value = 642
$ws.Range("D26").Value = 1566
$ws.Range("E26").Value = 211

# --- updated counts: Segovia (now row 27) ---------------------------------------------
$ws.Range("B27").Value = 2406
$ws.Range("C27").Value = 656
$ws.Range("D27").Value = 1578
$ws.Range("E27").Value = 172

# --- updated counts: Leon (now row 28) ------------------------------------------------
$ws.Range("B28").Value = 2403
$ws.Range("C28").Value = 1076
$ws.Range("D28").Value = 1024
$ws.Range("E28").Value = 303

# --- updated counts: Cantabria -----------------------------------------------------------
$ws.Range("B32").Value = 2160
$ws.Range("C32").Value = 718
$ws.Range("D32").Value = 1275
$ws.Range("E32").Value = 167

# --- updated counts: Gran Canaria ---------------------------------------------------------
$ws.Range("B33").Value = 2094
$ws.Range("C33").Value = 927
$ws.Range("D33").Value = 1046
$ws.Range("E33").Value = 121

# --- updated counts: Murcia ----------------------------------------------------------------
$ws.Range("B36").Value = 1695
$ws.Range("C36").Value = 761
$ws.Range("D36").Value = 811
$ws.Range("E36").Value = 123

# --- updated counts: Ceuta -------------------------------------------------------------------
$ws.Range("B58").Value = 118
$ws.Range("C58").Value = 76
$ws.Range("D58").Value = 38

# --- updated counts: Melilla -----------------------------------------------------------------
$ws.Range("C59").Value = 50
$ws.Range("D59").Value = 53
